$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dhh"
$ws.Range("C2").Value = "Ptch2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 2.179395
$ws.Range("H2").Value = 6.538185
$ws.Range("I2").Value = 0.4845018986408914
$ws.Range("J2").Value = 0.4845018986408914
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.217641
$ws.Range("N2").Value = 15.652923
$ws.Range("O2").Value = 0.3148443262844371
$ws.Range("P2").Value = 0.3148443262844371
$ws.Range("Q2").Value = 11.371300707195
$ws.Range("R2").Value = 102.341706364755
$ws.Range("S2").Value = 0.1525426738611221
$ws.Range("T2").Value = 0.1525426738611221

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dhh"
$ws.Range("C3").Value = "Ptch2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 2.179395
$ws.Range("H3").Value = 6.538185
$ws.Range("I3").Value = 0.4845018986408914
$ws.Range("J3").Value = 0.4845018986408914
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.194163666666666
$ws.Range("N3").Value = 27.582491
$ws.Range("O3").Value = 0.5547967492168427
$ws.Range("P3").Value = 0.5547967492168427
$ws.Range("Q3").Value = 20.037714324315
$ws.Range("R3").Value = 180.339428918835
$ws.Range("S3").Value = 0.2688000783553548
$ws.Range("T3").Value = 0.2688000783553548

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Dhh"
$ws.Range("C4").Value = "Ptch2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 2.179395
$ws.Range("H4").Value = 6.538185
$ws.Range("I4").Value = 0.4845018986408914
$ws.Range("J4").Value = 0.4845018986408914
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.160325
$ws.Range("N4").Value = 6.480975
$ws.Range("O4").Value = 0.1303589244987201
$ws.Range("P4").Value = 0.1303589244987201
$ws.Range("Q4").Value = 4.708201503374999
$ws.Range("R4").Value = 42.373813530375
$ws.Range("S4").Value = 0.06315914642441452
$ws.Range("T4").Value = 0.06315914642441452

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Dhh"
$ws.Range("C5").Value = "Ptch2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7348883333333333
$ws.Range("H5").Value = 2.204665
$ws.Range("I5").Value = 0.1633732264179005
$ws.Range("J5").Value = 0.1633732264179005
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.217641
$ws.Range("N5").Value = 15.652923
$ws.Range("O5").Value = 0.3148443262844371
$ws.Range("P5").Value = 0.3148443262844371
$ws.Range("Q5").Value = 3.834383498421666
$ws.Range("R5").Value = 34.509451485795
$ws.Range("S5").Value = 0.05143713340445868
$ws.Range("T5").Value = 0.05143713340445868

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Dhh"
$ws.Range("C6").Value = "Ptch2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7348883333333333
$ws.Range("H6").Value = 2.204665
$ws.Range("I6").Value = 0.1633732264179005
$ws.Range("J6").Value = 0.1633732264179005
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.194163666666666
$ws.Range("N6").Value = 27.582491
$ws.Range("O6").Value = 0.5547967492168427
$ws.Range("P6").Value = 0.5547967492168427
$ws.Range("Q6").Value = 6.756683613390555
$ws.Range("R6").Value = 60.810152520515
$ws.Range("S6").Value = 0.0906389349257184
$ws.Range("T6").Value = 0.0906389349257184

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Dhh"
$ws.Range("C7").Value = "Ptch2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7348883333333333
$ws.Range("H7").Value = 2.204665
$ws.Range("I7").Value = 0.1633732264179005
$ws.Range("J7").Value = 0.1633732264179005
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.160325
$ws.Range("N7").Value = 6.480975
$ws.Range("O7").Value = 0.1303589244987201
$ws.Range("P7").Value = 0.1303589244987201
$ws.Range("Q7").Value = 1.587597638708333
$ws.Range("R7").Value = 14.288378748375
$ws.Range("S7").Value = 0.0212971580877234
$ws.Range("T7").Value = 0.0212971580877234

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Dhh"
$ws.Range("C8").Value = "Ptch2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.583934333333333
$ws.Range("H8").Value = 4.751803
$ws.Range("I8").Value = 0.3521248749412083
$ws.Range("J8").Value = 0.3521248749412082
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.217641
$ws.Range("N8").Value = 15.652923
$ws.Range("O8").Value = 0.3148443262844371
$ws.Range("P8").Value = 0.3148443262844371
$ws.Range("Q8").Value = 8.264400718907666
$ws.Range("R8").Value = 74.379606470169
$ws.Range("S8").Value = 0.1108645190188564
$ws.Range("T8").Value = 0.1108645190188564

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Dhh"
$ws.Range("C9").Value = "Ptch2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.583934333333333
$ws.Range("H9").Value = 4.751803
$ws.Range("I9").Value = 0.3521248749412083
$ws.Range("J9").Value = 0.3521248749412082
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.194163666666666
$ws.Range("N9").Value = 27.582491
$ws.Range("O9").Value = 0.5547967492168427
$ws.Range("P9").Value = 0.5547967492168427
$ws.Range("Q9").Value = 14.56295149791922
$ws.Range("R9").Value = 131.066563481273
$ws.Range("S9").Value = 0.1953577359357696
$ws.Range("T9").Value = 0.1953577359357696

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Dhh"
$ws.Range("C10").Value = "Ptch2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.583934333333333
$ws.Range("H10").Value = 4.751803
$ws.Range("I10").Value = 0.3521248749412083
$ws.Range("J10").Value = 0.3521248749412082
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.160325
$ws.Range("N10").Value = 6.480975
$ws.Range("O10").Value = 0.1303589244987201
$ws.Range("P10").Value = 0.1303589244987201
$ws.Range("Q10").Value = 3.421812938658333
$ws.Range("R10").Value = 30.796316447925
$ws.Range("S10").Value = 0.04590261998658224
$ws.Range("T10").Value = 0.04590261998658223
